$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the default (unstyled) cell style used by all data rows so that
# forcing text storage on numeric-looking values does not change cell styling.
$origStyle = $ws.Range("D4").Style

$ws.Range('D2').Value = '56.518.74'
$ws.Range('E2').Value = '  -0.54%  '

$ws.Range('D3').Value = '2.319.16'
$ws.Range('E3').Value = '  -0.40%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '515.39'
$ws.Range('E5').Value = '  -1.20%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '131.68'
$ws.Range('E6').Value = '  -2.30%  '

$ws.Range('E7').Value = '  +0.41%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.532'
$ws.Range('E8').Value = '  -1.09%  '

$ws.Range('E9').Value = '  -2.68%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.153'
$ws.Range('E10').Value = '  -0.14%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.25'
$ws.Range('E11').Value = '  -1.32%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.336'
$ws.Range('E12').Value = '  -2.25%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '23.52'
$ws.Range('E13').Value = '  -1.87%  '

$ws.Range('D14').Value = '2.735.20'
$ws.Range('E14').Value = '  -0.39%  '

$ws.Range('D15').Value = '56.534.94'
$ws.Range('E15').Value = '  -0.62%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000132'
$ws.Range('E16').Value = '  -1.44%  '

$ws.Range('D17').Value = '2.324.99'
$ws.Range('E17').Value = '  -0.65%  '

$ws.Range('E18').Value = '  -1.63%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '327.42'
$ws.Range('E19').Value = '  +1.29%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.13'
$ws.Range('E20').Value = '  -2.34%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.73'
$ws.Range('E21').Value = '  +2.06%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.13%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '61.12'
$ws.Range('E23').Value = '  +0.53%  '

$ws.Range('E24').Value = '  +0.19%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.56'
$ws.Range('E25').Value = '  +7.36%  '

$ws.Range('E26').Value = '  +0.59%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.30'
$ws.Range('E27').Value = '  +0.04%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '167.71'
$ws.Range('E28').Value = '  +0.82%  '

$ws.Range('E29').Value = '  -2.16%  '

$ws.Range('D30').Value = '0.0₃0720'
$ws.Range('E30').Value = '  -3.13%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.13'
$ws.Range('E31').Value = '  -1.16%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '18.26'
$ws.Range('E32').Value = '  -0.65%  '

$ws.Range('E33').Value = '  -0.02%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  +0.70%  '

$ws.Range('E35').Value = '  -1.22%  '

$ws.Range('E36').Value = '  -2.80%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.882'
$ws.Range('E37').Value = '  -5.06%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.56'
$ws.Range('E38').Value = '  +0.12%  '

$ws.Range('E39').Value = '  +1.83%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '148.60'
$ws.Range('E40').Value = '  +7.12%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.373'
$ws.Range('E41').Value = '  -2.98%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.57'
$ws.Range('E42').Value = '  -1.06%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '279.08'
$ws.Range('E43').Value = '  -0.48%  '

$ws.Range('E44').Value = '  -6.10%  '

$ws.Range('E45').Value = '  -0.63%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0495'
$ws.Range('E46').Value = '  -2.34%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.557'
$ws.Range('E47').Value = '  -1.66%  '

$ws.Range('E48').Value = '  +1.84%  '

$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0214'
$ws.Range('E49').Value = '  -1.89%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '17.09'
$ws.Range('E50').Value = '  +0.69%  '

$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '11.01'
$ws.Range('E51').Value = '  +0.53%  '

# Restore the original (default) style on the Price column so the text-forcing
# number format above does not leave a stray style applied to these cells.
$ws.Range("D2:D51").Style = $origStyle
